$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "322.75"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.90%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.65"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.15%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.886"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "11.64%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08019"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.04%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.573"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.99%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.656"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.19%"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.922"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.26%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9315"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.30%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1220"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-9.91%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1957"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.25%"
$ws.Range("B12").Value = "MCDex"
$ws.Range("C12").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.720"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "19.65%"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09111"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-2.28%"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03540"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2.96%"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09565"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.15%"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001297"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-6.86%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006207"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2.00%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.353"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.26%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.950"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.29%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3560"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.96%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "8.03%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2409"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.29%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04410"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.66%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001262"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.39%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004411"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.16%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001141"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-11.58%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.12%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02429"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-3.32%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05225"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.47%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007419"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.53%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009378"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "8.78%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1404"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.89%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002121"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.79%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01122"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "37.79%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006743"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.33%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.09%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003003"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "3.58%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.09%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.09%"
